{"js": "// Load all paragraphs in the document body so we can inspect their text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three target paragraphs by their (unique) current text.\nlet para3 = null; // \"3. Cria\u00e7\u00e3o do boneco;\"\nlet para4 = null; // \"4. Dicas e palavra aleat\u00f3ria.\"\nlet para5 = null; // \"5. Mostrar em tela as letras reveladas.\"\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const t = p.text.trim();\n  if (t === \"3. Cria\u00e7\u00e3o do boneco;\") {\n    para3 = p;\n  } else if (t === \"4. Dicas e palavra aleat\u00f3ria.\") {\n    para4 = p;\n  } else if (t === \"5. Mostrar em tela as letras reveladas.\") {\n    para5 = p;\n  }\n}\n\nif (!para3 || !para4 || !para5) {\n  throw new Error(\"Could not locate one or more target paragraphs.\");\n}\n\n// 1) \"3. Cria\u00e7\u00e3o do boneco;\" -> \"3. Cria\u00e7\u00e3o do boneco [DONE];\"\npara3.insertText(\"3. Cria\u00e7\u00e3o do boneco [DONE];\", Word.InsertLocation.replace);\n\n// 2) \"4. Dicas e palavra aleat\u00f3ria.\" -> \"4. Dicas e palavra aleat\u00f3ria;\"\npara4.insertText(\"4. Dicas e palavra aleat\u00f3ria;\", Word.InsertLocation.replace);\n\n// 3) \"5. Mostrar em tela as letras reveladas.\" -> \"5. Mostrar em tela as letras reveladas [DONE].\"\npara5.insertText(\n  \"5. Mostrar em tela as letras reveladas [DONE].\",\n  Word.InsertLocation.replace\n);\n\n// 4) Insert a brand-new paragraph \"6. Refatorar o c\u00f3digo.\" right after paragraph 5,\n//    copying its run formatting (font, size).\nconst newPara = para5.insertParagraph(\"6. Refatorar o c\u00f3digo.\", Word.InsertLocation.after);\nnewPara.font.name = \"Helvetica\";\nnewPara.font.size = 29; // sz val 58 half-points = 29pt\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"3. Cria\u00e7\u00e3o do boneco;\" -> \"3. Cria\u00e7\u00e3o do boneco [DONE];\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"3. Cria\u00e7\u00e3o do boneco;\"\n$find.Replacement.Text = \"3. Cria\u00e7\u00e3o do boneco [DONE];\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n\n# 2) \"4. Dicas e palavra aleat\u00f3ria.\" -> \"4. Dicas e palavra aleat\u00f3ria;\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"4. Dicas e palavra aleat\u00f3ria.\"\n$find2.Replacement.Text = \"4. Dicas e palavra aleat\u00f3ria;\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2)\n\n# 3) \"5. Mostrar em tela as letras reveladas.\" -> \"5. Mostrar em tela as letras reveladas [DONE].\"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"5. Mostrar em tela as letras reveladas.\"\n$find3.Replacement.Text = \"5. Mostrar em tela as letras reveladas [DONE].\"\n$find3.Forward = $true\n$find3.Wrap = 1\n$find3.MatchCase = $true\n$find3.MatchWholeWord = $false\n$find3.Execute([ref]$find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find3.Replacement.Text, 2)\n\n# 4) Insert a brand-new paragraph \"6. Refatorar o c\u00f3digo.\" right after paragraph 5,\n#    matching the existing Helvetica/58 run formatting used throughout the list.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r`a\") -eq \"5. Mostrar em tela as letras reveladas [DONE].\") {\n        $rng = $p.Range\n        $rng.Collapse(0)  # wdCollapseEnd\n        $rng.InsertParagraphAfter()\n\n        $newPara = $p.Next()\n        $newRange = $newPara.Range\n        $newRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark itself\n        $newRange.Text = \"6. Refatorar o c\u00f3digo.\"\n\n        $newRange2 = $newPara.Range\n        $newRange2.MoveEnd(1, -1) | Out-Null\n        $newRange2.Font.Name = \"Helvetica\"\n        $newRange2.Font.NameBi = \"Helvetica\"\n        $newRange2.Font.Size = 29\n\n        break\n    }\n}\n"}
